# Update "想去人数" (F column) values on both the "展览" sheet and the
# "全部类型" sheet (which duplicates the same events), per the commit's
# refreshed scrape data.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (rId1 / sheet1.xml) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 1347
$wsExhibit.Range("F5").Value  = 1050
$wsExhibit.Range("F8").Value  = 1222
$wsExhibit.Range("F10").Value = 16
$wsExhibit.Range("F12").Value = 313
$wsExhibit.Range("F16").Value = 195
$wsExhibit.Range("F20").Value = 334
$wsExhibit.Range("F24").Value = 655
$wsExhibit.Range("F28").Value = 329
$wsExhibit.Range("F29").Value = 171
$wsExhibit.Range("F34").Value = 413

# --- Sheet: 全部类型 (rId4 / sheet4.xml) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 1347
$wsAll.Range("F6").Value  = 1050
$wsAll.Range("F9").Value  = 1222
$wsAll.Range("F12").Value = 16
$wsAll.Range("F14").Value = 313
$wsAll.Range("F18").Value = 195
$wsAll.Range("F25").Value = 334
$wsAll.Range("F32").Value = 655
$wsAll.Range("F36").Value = 329
$wsAll.Range("F39").Value = 171
$wsAll.Range("F48").Value = 413
